$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.431.53"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "3.443.24"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'574.54"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'143.95"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D7").Value = "3.443.46"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.479"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "4.034.43"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "'28.45"
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "3.442.82"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "61.555.81"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  +3.61%  "
$ws.Range("D20").Value = "'14.39"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "'9.37"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "'397.99"
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'74.02"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").Value = "'0.996"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'0.0000122"
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("D27").Value = "3.583.22"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "'0.178"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "'7.64"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'1.48"
$ws.Range("E31").Value = "  -7.24%  "
$ws.Range("D32").Value = "'8.21"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'23.96"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.473.90"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'7.03"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'5.15"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").Value = "'167.13"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'0.0789"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "'27.90"
$ws.Range("E42").Value = "  +6.06%  "
$ws.Range("D43").Value = "'0.802"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "2.623.23"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").Value = "'23.05"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("E51").Value = "  +2.68%  "
